# Auto-generated edit script applying numeric corrections described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1999.5
$ws.Range("J7").Value = 1999.5
$ws.Range("L7").Value = 1999.5
$ws.Range("N7").Value = -2223.5

$ws.Range("H11").Value = 40000320
$ws.Range("I11").Value = 40000320
$ws.Range("K11").Value = 40000320
$ws.Range("M11").Value = -40000180

$ws.Range("H14").Value = 1999.5
$ws.Range("J14").Value = 1999.5
$ws.Range("L14").Value = 1999.5
$ws.Range("N14").Value = -2381.5

$ws.Range("H40").Value = 2469.8
$ws.Range("J40").Value = 2116.6667
$ws.Range("L40").Value = 2116.6667
$ws.Range("N40").Value = -2466.6667

$ws.Range("H80").Value = 15798.429
$ws.Range("I80").Value = 2157
$ws.Range("J80").Value = 33987
$ws.Range("K80").Value = 6471
$ws.Range("L80").Value = 101961
$ws.Range("M80").Value = -5473
$ws.Range("N80").Value = -103957

$ws.Range("H83").Value = 15798.429
$ws.Range("I83").Value = 2157
$ws.Range("J83").Value = 33987
$ws.Range("K83").Value = 19413
$ws.Range("L83").Value = 305883
$ws.Range("M83").Value = -14421
$ws.Range("N83").Value = -315867

$ws.Range("H92").Value = 947396.4
$ws.Range("I92").Value = 1231414.6
$ws.Range("K92").Value = 1231414.6
$ws.Range("M92").Value = -1230166.6

$ws.Range("H99").Value = 777.6667
$ws.Range("I99").Value = 345.25
$ws.Range("J99").Value = 1642.5
$ws.Range("K99").Value = 1035.75
$ws.Range("L99").Value = 4927.5
$ws.Range("M99").Value = 462.25
$ws.Range("N99").Value = -7923.5

$ws.Range("H129").Value = 900.55554
$ws.Range("J129").Value = 875.58826
$ws.Range("L129").Value = 2626.76478
$ws.Range("N129").Value = -12626.76478

$ws.Range("H131").Value = 3369.0833
$ws.Range("I131").Value = 923.8
$ws.Range("J131").Value = 5115.7144
$ws.Range("K131").Value = 2771.4
$ws.Range("L131").Value = 15347.1432
$ws.Range("M131").Value = 2268.6
$ws.Range("N131").Value = -25427.1432

$ws.Range("H137").Value = 3531.3333
$ws.Range("I137").Value = 1370.125
$ws.Range("K137").Value = 4110.375
$ws.Range("M137").Value = -1560.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 80000
$ws.Range("J24").Value = 80000
$ws.Range("L24").Value = 80000
$ws.Range("N24").Value = -80748

$ws.Range("H61").Value = 5150.0557
$ws.Range("J61").Value = 10937.4
$ws.Range("L61").Value = 10937.4
$ws.Range("N61").Value = -11361.4

$ws.Range("H100").Value = 80000
$ws.Range("J100").Value = 80000
$ws.Range("L100").Value = 80000
$ws.Range("N100").Value = -82164

$ws.Range("H122").Value = 2418.5557
$ws.Range("I122").Value = 2375.5
$ws.Range("K122").Value = 7126.5
$ws.Range("M122").Value = -4676.5

$ws.Range("H136").Value = 5150.0557
$ws.Range("J136").Value = 10937.4
$ws.Range("L136").Value = 32812.2
$ws.Range("N136").Value = -37912.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1998.1666
$ws.Range("I105").Value = 2074.647
$ws.Range("J105").Value = 698
$ws.Range("K105").Value = 2074.647
$ws.Range("L105").Value = 698
$ws.Range("M105").Value = -327.6469999999999
$ws.Range("N105").Value = -4192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5268.8887
$ws.Range("J3").Value = 8680
$ws.Range("L3").Value = 8680
$ws.Range("N3").Value = -8906

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H31").Value = 2666.2646
$ws.Range("I31").Value = 1896.8462
$ws.Range("K31").Value = 1896.8462
$ws.Range("M31").Value = -1601.8462

$ws.Range("H34").Value = 2666.2646
$ws.Range("I34").Value = 1896.8462
$ws.Range("K34").Value = 1896.8462
$ws.Range("M34").Value = -1694.8462

$ws.Range("H74").Value = 29333.334
$ws.Range("J74").Value = 29333.334
$ws.Range("L74").Value = 29333.334
$ws.Range("N74").Value = -31081.334

$ws.Range("H77").Value = 29333.334
$ws.Range("J77").Value = 29333.334
$ws.Range("L77").Value = 88000.00199999999
$ws.Range("N77").Value = -96736.00199999999

$ws.Range("H122").Value = 1355.4375
$ws.Range("I122").Value = 1307.9524
$ws.Range("J122").Value = 1446.091
$ws.Range("K122").Value = 3923.857199999999
$ws.Range("L122").Value = 4338.272999999999
$ws.Range("M122").Value = -1473.857199999999
$ws.Range("N122").Value = -9238.272999999999

$ws.Range("H134").Value = 851.0714
$ws.Range("I134").Value = 839.0769
$ws.Range("J134").Value = 1007
$ws.Range("K134").Value = 2517.2307
$ws.Range("L134").Value = 3021
$ws.Range("M134").Value = 17.76929999999993
$ws.Range("N134").Value = -8091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12297.561
$ws.Range("I4").Value = 12297.561
$ws.Range("K4").Value = 36892.683
$ws.Range("M4").Value = -36780.683

$ws.Range("H36").Value = 1674.5
$ws.Range("I36").Value = 1674.5
$ws.Range("K36").Value = 5023.5
$ws.Range("M36").Value = -4854.5

$ws.Range("H56").Value = 7486.846
$ws.Range("I56").Value = 7486.846
$ws.Range("K56").Value = 7486.846
$ws.Range("M56").Value = -6956.846

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws.Range("H131").Value = 785.97
$ws.Range("J131").Value = 797.9474
$ws.Range("L131").Value = 2393.8422
$ws.Range("N131").Value = -12473.8422

$ws.Range("H133").Value = 2282.5
$ws.Range("I133").Value = 1376.6666
$ws.Range("K133").Value = 4129.9998
$ws.Range("M133").Value = 930.0002000000004

$ws.Range("H139").Value = 22031.6
$ws.Range("I139").Value = 27040
$ws.Range("K139").Value = 81120
$ws.Range("M139").Value = -75980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9125
$ws.Range("J5").Value = 9125
$ws.Range("L5").Value = 9125
$ws.Range("N5").Value = -9349

$ws.Range("H126").Value = 2695730
$ws.Range("I126").Value = 4276154.5
$ws.Range("K126").Value = 12828463.5
$ws.Range("M126").Value = -12825993.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 291350
$ws.Range("J2").Value = 82700
$ws.Range("L2").Value = 82700
$ws.Range("N2").Value = -82924

$ws.Range("H122").Value = 5381.4375
$ws.Range("J122").Value = 6288.8887
$ws.Range("L122").Value = 18866.6661
$ws.Range("N122").Value = -23766.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1162.75
$ws.Range("I113").Value = 943.2857
$ws.Range("J113").Value = 1470
$ws.Range("K113").Value = 2829.8571
$ws.Range("L113").Value = 4410
$ws.Range("M113").Value = -659.8571000000002
$ws.Range("N113").Value = -8750

$ws.Range("H122").Value = 111760.47
$ws.Range("I122").Value = 118589.25
$ws.Range("K122").Value = 355767.75
$ws.Range("M122").Value = -353317.75
